$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E2").Value = 69
$ws.Range("F2").Value = 47
$ws.Range("H2").Value = 59
$ws.Range("E3").Value = 30
$ws.Range("E7").Value = 17
$ws.Range("F7").Value = 6
$ws.Range("H7").Value = 7
$ws.Range("E8").Value = 19
$ws.Range("F8").Value = 11
$ws.Range("H8").Value = 12
$ws.Range("E12").Value = 12
$ws.Range("F12").Value = 6
$ws.Range("H12").Value = 6
$ws.Range("E15").Value = 189
$ws.Range("F15").Value = 111
$ws.Range("H15").Value = 152
$ws.Range("E17").Value = 146
$ws.Range("F17").Value = 78
$ws.Range("H17").Value = 110
$ws.Range("E18").Value = 141
$ws.Range("F18").Value = 71
$ws.Range("H18").Value = 108
$ws.Range("E19").Value = 75
$ws.Range("F19").Value = 47
$ws.Range("H19").Value = 60
$ws.Range("E20").Value = 8
$ws.Range("F20").Value = 5
$ws.Range("H20").Value = 8
$ws.Range("E25").Value = 30
$ws.Range("F25").Value = 19
$ws.Range("H25").Value = 27
$ws.Range("E26").Value = 38
$ws.Range("E27").Value = 20
$ws.Range("F27").Value = 16
$ws.Range("H27").Value = 20
$ws.Range("E28").Value = 25
$ws.Range("E33").Value = 49
$ws.Range("E34").Value = 29
$ws.Range("F34").Value = 12
$ws.Range("H34").Value = 15
$ws.Range("E35").Value = 15
$ws.Range("F35").Value = 9
$ws.Range("H35").Value = 10
$ws.Range("E36").Value = 130
$ws.Range("F36").Value = 67
$ws.Range("H36").Value = 99
$ws.Range("E37").Value = 66
$ws.Range("F37").Value = 42
$ws.Range("H37").Value = 54
$ws.Range("E38").Value = 94
$ws.Range("E43").Value = 35
$ws.Range("E45").Value = 30
$ws.Range("F45").Value = 16
$ws.Range("H45").Value = 23
$ws.Range("E47").Value = 70
$ws.Range("F47").Value = 46
$ws.Range("H47").Value = 56
$ws.Range("F62").Value = 19
$ws.Range("H62").Value = 33
$ws.Range("E63").Value = 51
$ws.Range("F63").Value = 22
$ws.Range("H63").Value = 30
$ws.Range("E64").Value = 43
$ws.Range("F64").Value = 22
$ws.Range("H64").Value = 27
$ws.Range("E66").Value = 42
$ws.Range("E67").Value = 47
$ws.Range("F67").Value = 27
$ws.Range("H67").Value = 35
$ws.Range("E69").Value = 19
$ws.Range("E71").Value = 49
$ws.Range("E72").Value = 54
$ws.Range("E73").Value = 38
$ws.Range("F73").Value = 18
$ws.Range("H73").Value = 30
$ws.Range("E77").Value = 68
$ws.Range("F77").Value = 28
$ws.Range("H77").Value = 45
$ws.Range("E78").Value = 51
$ws.Range("E79").Value = 49
$ws.Range("F81").Value = 15
$ws.Range("H81").Value = 20
$ws.Range("E83").Value = 14
$ws.Range("E88").Value = 39
